$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow a couple of columns (J and AB) from width 8 to width 7 ---
# Excel's ColumnWidth property is offset from the raw OOXML <col width> by
# ~0.8333333333333334 (5/6) character units in this runtime, so compensate.
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667
$ws.Columns.Item(28).ColumnWidth = 6.166666666666667

# --- Round row 5's measurement values to 2 decimal places ("custom accuracy") ---
$ws.Range("C5").Value = 16.98
$ws.Range("D5").Value = 1.33
$ws.Range("E5").Value = 50.12
$ws.Range("F5").Value = 40.95
$ws.Range("G5").Value = 18.15
$ws.Range("H5").Value = 68.22
$ws.Range("I5").Value = 27.92
$ws.Range("J5").Value = 12.38
$ws.Range("L5").Value = 20.11
$ws.Range("M5").Value = 21.2
$ws.Range("N5").Value = 5.8
$ws.Range("O5").Value = 18.05
$ws.Range("P5").Value = 25.66
$ws.Range("Q5").Value = 15.22
$ws.Range("R5").Value = 0.86
$ws.Range("S5").Value = 0.92
$ws.Range("T5").Value = 267.37
$ws.Range("U5").Value = 50.42
$ws.Range("V5").Value = 16.66
$ws.Range("W5").Value = 33.87
$ws.Range("X5").Value = 17.75
$ws.Range("Y5").Value = 2.69
$ws.Range("Z5").Value = 33.5
$ws.Range("AA5").Value = 14.71
$ws.Range("AB5").Value = 13.06
$ws.Range("AC5").Value = 15.35
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 61.91
$ws.Range("AG5").Value = 9.37
$ws.Range("AH5").Value = 20.83

# --- Remove the last data row (row 6); the dataset now ends at row 5 ---
$ws.Rows.Item(6).Delete()
